$wb = $excel.ActiveWorkbook

# --- Update the "Activity" sheet: rename the Company-Discussed meeting
#     test data to refer to "Engagement" instead of "Company", and move
#     the sheet's selection/active state here. ---
$activity = $wb.Worksheets.Item("Activity")
$activity.Range("B2").Value = "Engagement Discussed Meeting"
$activity.Range("E2").Value = "Engagement Discussed Meeting Description"

# --- Switch the active/selected sheet from "AddContact" to "Activity" ---
$addContact = $wb.Worksheets.Item("AddContact")
[void]$addContact.Range("J14").Select()

[void]$activity.Activate()
[void]$activity.Range("E11").Select()
